# Add the "metadata" worksheet after the existing "data" sheet.
$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item(1)
$meta = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $dataSheet)
$meta.Name = "metadata"

# match the page-margin defaults (inches) openpyxl used on the "data" sheet
$meta.PageSetup.LeftMargin = 54
$meta.PageSetup.RightMargin = 54
$meta.PageSetup.TopMargin = 72
$meta.PageSetup.BottomMargin = 72
$meta.PageSetup.HeaderMargin = 36
$meta.PageSetup.FooterMargin = 36

# ---- header row ----
$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

# ---- data row ----
$meta.Range("A2").Value = 0
$meta.Range("B2").Value = "Gastrointestinal neuromuscular disorders"
$meta.Range("C2").Value = 61
# data_version ("1.15") is stored as TEXT, not a number - build it as a
# formula that yields a string, then flatten to a plain value so no
# residual number-format style gets left behind on the cell.
$meta.Range("D2").Formula = "=""1.15"""
$meta.Range("D2").Copy() | Out-Null
$meta.Range("D2").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = 0
$meta.Range("E2").Value = "2021-07-12T10:44:17.164775Z"
$meta.Range("F2").Value = "2021-10-05 14:20:27.764362"
$meta.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/61/?format=json"

# ---- copy the bold/centered/bordered header style from the "data" sheet ----
$dataSheet.Range("B1").Copy() | Out-Null
$meta.Range("B1:G1").PasteSpecial(-4122) | Out-Null
$dataSheet.Range("A2").Copy() | Out-Null
$meta.Range("A2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ---- refresh the "time_taken" timestamps on the data sheet (re-queried values) ----
$dataSheet.Range("F2").Value = "2021-10-05 14:20:27.768013"
$dataSheet.Range("F3").Value = "2021-10-05 14:20:27.768021"
$dataSheet.Range("F4").Value = "2021-10-05 14:20:27.768024"
$dataSheet.Range("F5").Value = "2021-10-05 14:20:27.768027"
$dataSheet.Range("F6").Value = "2021-10-05 14:20:27.768030"
$dataSheet.Range("F7").Value = "2021-10-05 14:20:27.768032"
$dataSheet.Range("F8").Value = "2021-10-05 14:20:27.768035"
$dataSheet.Range("F9").Value = "2021-10-05 14:20:27.768037"
$dataSheet.Range("F10").Value = "2021-10-05 14:20:27.768040"
$dataSheet.Range("F11").Value = "2021-10-05 14:20:27.768042"
$dataSheet.Range("F12").Value = "2021-10-05 14:20:27.768045"
$dataSheet.Range("F13").Value = "2021-10-05 14:20:27.768047"
$dataSheet.Range("F14").Value = "2021-10-05 14:20:27.768050"
$dataSheet.Range("F15").Value = "2021-10-05 14:20:27.768052"
$dataSheet.Range("F16").Value = "2021-10-05 14:20:27.768055"
$dataSheet.Range("F17").Value = "2021-10-05 14:20:27.768057"
$dataSheet.Range("F18").Value = "2021-10-05 14:20:27.768060"
$dataSheet.Range("F19").Value = "2021-10-05 14:20:27.768062"
$dataSheet.Range("F20").Value = "2021-10-05 14:20:27.768065"
$dataSheet.Range("F21").Value = "2021-10-05 14:20:27.768067"
$dataSheet.Range("F22").Value = "2021-10-05 14:20:27.768070"
$dataSheet.Range("F23").Value = "2021-10-05 14:20:27.768072"
$dataSheet.Range("F24").Value = "2021-10-05 14:20:27.768074"
$dataSheet.Range("F25").Value = "2021-10-05 14:20:27.768077"
$dataSheet.Range("F26").Value = "2021-10-05 14:20:27.768080"
$dataSheet.Range("F27").Value = "2021-10-05 14:20:27.768082"
$dataSheet.Range("F28").Value = "2021-10-05 14:20:27.768085"
$dataSheet.Range("F29").Value = "2021-10-05 14:20:27.768087"
$dataSheet.Range("F30").Value = "2021-10-05 14:20:27.768089"
$dataSheet.Range("F31").Value = "2021-10-05 14:20:27.768092"

# ---- make "data" the active sheet (matches original workbook view) ----
$dataSheet.Select()
